# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.143.28"
$ws.Range("E2").Value = "  +2.38%  "

$ws.Range("D3").Value = "3.400.15"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'559.96"
$ws.Range("E5").Value = "  +2.73%  "

$ws.Range("D6").Value = "'175.60"
$ws.Range("E6").Value = "  +2.38%  "

$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +2.42%  "

$ws.Range("D8").Value = "3.392.31"
$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("D9").Value = "'1.00"

$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  +12.26%  "

$ws.Range("D11").Value = "'0.631"
$ws.Range("E11").Value = "  +3.69%  "

$ws.Range("D12").Value = "'54.91"
$ws.Range("E12").Value = "  +2.97%  "

$ws.Range("D13").Value = "'0.0000280"
$ws.Range("E13").Value = "  +5.86%  "

$ws.Range("D14").Value = "'9.15"
$ws.Range("E14").Value = "  +3.38%  "

$ws.Range("D15").Value = "3.946.19"
$ws.Range("E15").Value = "  +8.56%  "

$ws.Range("D16").Value = "'18.40"
$ws.Range("E16").Value = "  +2.32%  "

$ws.Range("D17").Value = "3.413.64"
$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("E18").Value = "  +2.16%  "

$ws.Range("D19").Value = "65.114.65"
$ws.Range("E19").Value = "  +2.45%  "

$ws.Range("D20").Value = "'11.87"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").Value = "'0.995"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").Value = "'474.02"
$ws.Range("E22").Value = "  +15.84%  "

$ws.Range("D23").Value = "'4.97"
$ws.Range("E23").Value = "  +14.13%  "

$ws.Range("D24").Value = "'4.14"
$ws.Range("E24").Value = "  +3.25%  "

$ws.Range("D25").Value = "'87.21"
$ws.Range("E25").Value = "  +5.27%  "

$ws.Range("D26").Value = "'13.46"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "  +7.45%  "

$ws.Range("D28").Value = "'10.89"
$ws.Range("E28").Value = "  +3.73%  "

$ws.Range("D29").Value = "'8.82"
$ws.Range("E29").Value = "  +3.09%  "

$ws.Range("D30").Value = "'31.24"
$ws.Range("E30").Value = "  +7.82%  "

$ws.Range("D31").Value = "'6.71"
$ws.Range("E31").Value = "  +5.93%  "

$ws.Range("D32").Value = "'11.55"
$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("D33").Value = "'61.86"
$ws.Range("E33").Value = "  +7.31%  "

$ws.Range("D34").Value = "'573.43"
$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("D35").Value = "'0.108"
$ws.Range("E35").Value = "  +2.40%  "

$ws.Range("D36").Value = "'0.999"

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.140"
$ws.Range("E37").Value = "  -4.61%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.53"
$ws.Range("E38").Value = "  +3.92%  "

$ws.Range("D39").Value = "0.0₃0763"
$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("D40").Value = "'35.85"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("D41").Value = "'0.372"
$ws.Range("E41").Value = "  +2.10%  "

$ws.Range("D42").Value = "3.098.09"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  +3.15%  "

$ws.Range("D45").Value = "'0.0416"
$ws.Range("E45").Value = "  +4.49%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.136"
$ws.Range("E46").Value = "  +6.46%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.48"
$ws.Range("E47").Value = "  +3.03%  "

$ws.Range("D48").Value = "'3.16"
$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("D50").Value = "'137.57"
$ws.Range("E50").Value = "  +3.88%  "

$ws.Range("D51").Value = "'8.33"
$ws.Range("E51").Value = "  +3.82%  "
